# Fruta / hortaliza, semanal
#
# The source data table (row 2..910) gets one stale record (old row 859)
# dropped and five fresh weekly records inserted in its place, pushing
# everything from the old row 860 onward down by four rows (910 -> 914).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the stale record that used to sit at row 859.
$ws.Rows.Item(859).Delete()

# 2) Make room for the five new records right above the old row 860
#    (which, after the delete above, is sitting at row 859).
$ws.Range("A859:A863").EntireRow.Insert()

# 3) Populate the five new records (columns A-T), keeping the same
#    Mercado/Producto/Categoria metadata used throughout this block.
$newRows = @(
    @{ Row=859; D=45267; L="Especial"; M=100; N=12000; O=12000; P=12000; Q="$/bandeja 7 kilos"; R="Provincia de Melipilla"; S=1714 },
    @{ Row=860; D=45267; L="Especial"; M=80;  N=12000; O=12000; P=12000; Q="$/caja 7 kilos";    R="Región del Maule";       S=1714 },
    @{ Row=861; D=45267; L="Primera";  M=200; N=10000; O=10000; P=10000; Q="$/bandeja 7 kilos"; R="Provincia de Melipilla"; S=1429 },
    @{ Row=862; D=45267; L="Primera";  M=120; N=10000; O=10000; P=10000; Q="$/caja 7 kilos";    R="Región del Maule";       S=1429 },
    @{ Row=863; D=45267; L="Especial"; M=120; N=9000;  O=9000;  P=9000;  Q="$/bandeja 7 kilos"; R="Provincia de Melipilla"; S=1286 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value = "Maule"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 7
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 7
}
